# "-be armor" (under Accessories) is done -> strike it through, the same
# way the other finished goals in this list are marked.
$d = $word.ActiveDocument

# The bullet is currently split across two runs ("-" and "be armor"). Find
# the whole phrase, delete it, and retype it in one shot so it collapses
# back down into a single run before we apply the strikethrough formatting.
$target = $d.Content
$found = $target.Find.Execute("-be armor", $false, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $target.Paragraphs(1)
    $target.Delete()
    $para.Range.InsertBefore("-be armor")
}

# Re-locate the now-merged run and apply the strikethrough (keeping the
# Calibri font it already had).
$line = $d.Content
$line.Find.Execute("-be armor", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$line.Font.Name = "Calibri"
$line.Font.StrikeThrough = 1

# Strike the paragraph mark too, so the whole line -- not just the visible
# characters -- carries the struck-through formatting.
$line.Paragraphs(1).Range.Font.StrikeThrough = 1
